# Edit script: update report-generated timestamp and zero-out billed
# amount / pricing totals (reflecting a no-violation / no-billable-work
# scenario for this week's report).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Work Report")

# D5: Update the "Report Generated On" timestamp.
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:49 PM"

# C8: Total Billed Amount -> 0
$ws.Range("C8").Value = 0

# H16: Line item pricing -> 0
$ws.Range("H16").Value = 0

# H17: TOTAL pricing -> 0
$ws.Range("H17").Value = 0
